# Saldo_guide.xlsx update
# - Roll the "Dt. Referencia" (column G) reference date forward one day
#   (2024-09-10 -> 2024-09-11) for every data row.
# - Refresh the three figures that actually moved for account 376732 /
#   2876742152 / 35132728120 (rows 58, 103, 104).
# - Rename the sheet (and tab) to match the new extraction timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All data rows run from row 2 to row 274; column G holds the reference
# date as the Excel serial 45545 (2024-09-10). Bump every one of them to
# 45546 (2024-09-11) in a single range write.
$ws.Range("G2:G274").Value = 45546

# Row 58: Saldo Previsto / Vl. Total increased to 6549.26
$ws.Cells.Item(58, 5).Value = 6549.26
$ws.Cells.Item(58, 8).Value = 6549.26

# Row 103: Saldo Previsto / Vl. Total increased to 2280.2399999999998
$ws.Cells.Item(103, 5).Value = 2280.2399999999998
$ws.Cells.Item(103, 8).Value = 2280.2399999999998

# Row 104: Vl. Projetado dropped to 0, Saldo Previsto recalculated to
# 4780.63 (Vl. Total, column H, was already 4780.63 and is unchanged)
$ws.Cells.Item(104, 4).Value = 0
$ws.Cells.Item(104, 5).Value = 4780.63

# Rename sheet/tab to the new export run id
$ws.Name = "IClientBalance-20240911-095805-"
